$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1, matching the style used by the other header cells (e.g. J1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# New data values for column K, rows 2-8 (plain/unstyled, like the other data cells)
$values = @("PROCEDURE", "DEVICE", "PROCEDURE", "OTHER", "OTHER", "DRUG", "PROCEDURE")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
